# Yii2RbacTutorial.pptx - "adding merit of RBAC" edit
# Slide 11 ("Why we want to use a framework like RBAC") contains two
# code-sample text boxes. This edit:
#   1. Lower-cases "If (" -> "if (" in both code samples.
#   2. Grows the second code sample's text box to fit a new bullet.
#   3. Appends a new bold bullet paragraph to the second code sample.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# Shape 2: first code sample - "If ($post->created_by == ...) {"
$sh1 = $s.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$full1 = $tr1.Text
$needle1 = "If (`$post->"
$pos1 = $full1.IndexOf($needle1)
$run1 = $tr1.Characters($pos1 + 1, $needle1.Length)
$run1.Text = "if (`$post->"

# Shape 3: second code sample - "If (Yii::$app->user->can('updatePost')) {"
$sh2 = $s.Shapes.Item(3)
$tr2 = $sh2.TextFrame.TextRange
$full2 = $tr2.Text
$needle2 = "If ("
$pos2 = $full2.IndexOf($needle2)
$run2 = $tr2.Characters($pos2 + 1, $needle2.Length)
$run2.Text = "if ("

# Grow the text box to make room for the new bullet paragraph.
$sh2.Height = 138.1359

# Append the new bullet as its own paragraph at the end of the text box.
$lastParaIndex = $tr2.Paragraphs().Count
$lastPara = $tr2.Paragraphs($lastParaIndex, 1)
$rsquo = [char]0x2019
$hellip = [char]0x2026
$newBullet = [char]13 + "+ It" + $rsquo + "s easier to maintain (adding more privileges, role" + $hellip + ")"
$newRun = $lastPara.InsertAfter($newBullet)
